# Applies the two changes described by the commit "updated figs and ms":
#   1. The cached text of every "datetimeFigureOut" date field (Slide
#      Master, every Slide Layout, and the Notes Master) is bumped from
#      15/11/2022 -> 11/12/2022.
#   2. The caption textbox on slide 1 is reworded from
#      "Temporal variability of species detections" to
#      "Temporal stability of species detections".

$p = $ppt.ActivePresentation

$oldDate = "15/11/2022"
$newDate = "11/12/2022"

# ppPlaceholderDate
$ppPlaceholderDate = 16

function Update-DatePlaceholder {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }

        if ($isDatePlaceholder -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# 1a. Slide Master date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# 1b. Every Slide Layout's date placeholder.
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# 1c. Notes Master date placeholder.
$notesMaster = $p.NotesMaster
Update-DatePlaceholder $notesMaster.Shapes

# 2. Slide 1 caption textbox wording tweak.
$oldCaption = "Temporal variability of species detections"
$newCaption = "Temporal stability of species detections"

$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq $oldCaption) {
            $tr.Text = $newCaption
        }
    }
}
